# Fill in the new sex (U = unrecorded/unknown), sire-age (F) and dam-age (G)
# columns for each data row, and update the selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  F = 1;  G = 1 }
    @{ Row = 3;  F = 3;  G = 1 }
    @{ Row = 4;  F = 5;  G = 1 }
    @{ Row = 5;  F = 3;  G = 5 }
    @{ Row = 6;  F = 5;  G = 5 }
    @{ Row = 7;  F = 6;  G = 6 }
    @{ Row = 8;  F = 8;  G = 6 }
    @{ Row = 9;  F = 10; G = 6 }
    @{ Row = 10; F = 8;  G = 10 }
    @{ Row = 11; F = 10; G = 10 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = "U"
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}

[void]$ws.Range("D2:G11").Select()
